$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the shared "block" values (B, C, I columns) for every data row (7-11) ---
# CR-PK-CUS-POC-2310038 -> CR-PK-CUS-POC-2310041
$ws.Range("B7:B11").Value = "CR-PK-CUS-POC-2310041"
# s1702 -> s1704
$ws.Range("C7:C11").Value = "s1704"
# PKSUPTOPKCUS17001-s1702-002 -> PKSUPTOPKCUS17001-s1704-004
$ws.Range("I7:I11").Value = "PKSUPTOPKCUS17001-s1704-004"

# --- Re-sequence the per-row scenario blocks (E/F/K columns) ---
# Row 7 was scenario "003" -> becomes scenario "001"
$ws.Range("E7").Value = "scenario1720230604001"
$ws.Range("F7").Value = "PK-CUS-scenario17-20230604-001"
$ws.Range("K7").Value = "PK-SUP-scenario17-20230604-001"

# Row 8 was scenario "002" -> becomes scenario "003"
$ws.Range("E8").Value = "scenario1720230604003"
$ws.Range("F8").Value = "PK-CUS-scenario17-20230604-003"
$ws.Range("K8").Value = "PK-SUP-scenario17-20230604-003"

# Row 9 was scenario "005" -> becomes scenario "002"
$ws.Range("E9").Value = "scenario1720230604002"
$ws.Range("F9").Value = "PK-CUS-scenario17-20230604-002"
$ws.Range("K9").Value = "PK-SUP-scenario17-20230604-002"

# Row 10 was scenario "004" -> becomes scenario "005"
$ws.Range("E10").Value = "scenario1720230604005"
$ws.Range("F10").Value = "PK-CUS-scenario17-20230604-005"
$ws.Range("K10").Value = "PK-SUP-scenario17-20230604-005"

# Row 11 was scenario "001" -> becomes scenario "004"
$ws.Range("E11").Value = "scenario1720230604004"
$ws.Range("F11").Value = "PK-CUS-scenario17-20230604-004"
$ws.Range("K11").Value = "PK-SUP-scenario17-20230604-004"

# --- Resolve the H9/H10 "Back No." conflict: the two cells swap their blank-vs-
#     empty-text storage. H9 goes from a true blank to a stored empty-text value;
#     H10 goes from a stored empty-text value to a true blank. ---

# H9: true blank -> stored empty string (keep formatting identical to its
# neighbours by re-applying H8's format after the text-value assignment, since
# a bare "" assignment is interpreted as ClearContents rather than empty text).
$ws.Range("H9").Value = "'"
$ws.Range("H8").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# H10: stored empty string -> true blank
$ws.Range("H10").ClearContents()
